# Vega Modelo de Temuco - Jengibre: add a new weekly price record.
# A new row is inserted at row 124, pushing the existing rows 124-144
# down to 125-145, and the new row 124 is filled with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(124).Insert()

$ws.Range("A124").Value = 10
$ws.Range("B124").Value = "Vega Modelo de Temuco"
$ws.Range("C124").Value = "La Araucanía"
$ws.Range("D124").Value = 44637
$ws.Range("E124").Value = 9
$ws.Range("F124").Value = 100114007
$ws.Range("G124").Value = "Jengibre"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 70
$ws.Range("K124").Value = 22000
$ws.Range("L124").Value = 25000
$ws.Range("M124").Value = 23286
$ws.Range("N124").Value = "$/caja 13 kilos"
$ws.Range("O124").Value = "Perú"
$ws.Range("P124").Value = 1791
$ws.Range("Q124").Value = 13
$ws.Range("R124").Value = "Hortaliza"
